# Continuation of work on the "eb" project: rename the worksheet to the
# generic name used going forward, and start wiring up iterative
# calculation (set the max-change tolerance) without turning iteration on
# yet -- matching the "Not finished" state described in the commit.

$wb = $excel.ActiveWorkbook

# Rename the only worksheet from "cond_eb1_c" to "cond".
$wb.Worksheets.Item(1).Name = "cond"

# Started configuring iterative calculation: set the maximum change
# tolerance (Formulas > Maximum Change) to 0.0001, while leaving the
# "Enable iterative calculation" checkbox itself untouched/off for now.
$excel.MaxChange = 0.0001
